# Fruta / hortaliza, semanal
# Insert two new weekly price rows for "Murcott" (Primera / Segunda) at the
# top of the Vega Modelo de Temuco block, pushing the existing rows
# (349-367) down to (351-369).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 349 (old data shifts down).
$ws.Rows.Item(349).Insert()
$ws.Rows.Item(349).Insert()

# ---- New row 349: Murcott / Primera ----
$ws.Cells.Item(349, 1).Value = 10
$ws.Cells.Item(349, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(349, 3).Value = "La Araucanía"
$ws.Cells.Item(349, 4).Value = 44461
$ws.Cells.Item(349, 5).Value = 9
$ws.Cells.Item(349, 6).Value = "Fruta"
$ws.Cells.Item(349, 7).Value = 100102
$ws.Cells.Item(349, 8).Value = "Cítricos"
$ws.Cells.Item(349, 9).Value = 100102004
$ws.Cells.Item(349, 10).Value = "Mandarina"
$ws.Cells.Item(349, 11).Value = "Murcott"
$ws.Cells.Item(349, 12).Value = "Primera"
$ws.Cells.Item(349, 13).Value = 200
$ws.Cells.Item(349, 14).Value = 12000
$ws.Cells.Item(349, 15).Value = 12000
$ws.Cells.Item(349, 16).Value = 12000
$ws.Cells.Item(349, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(349, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(349, 19).Value = 667
$ws.Cells.Item(349, 20).Value = 18

# ---- New row 350: Murcott / Segunda ----
$ws.Cells.Item(350, 1).Value = 10
$ws.Cells.Item(350, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(350, 3).Value = "La Araucanía"
$ws.Cells.Item(350, 4).Value = 44461
$ws.Cells.Item(350, 5).Value = 9
$ws.Cells.Item(350, 6).Value = "Fruta"
$ws.Cells.Item(350, 7).Value = 100102
$ws.Cells.Item(350, 8).Value = "Cítricos"
$ws.Cells.Item(350, 9).Value = 100102004
$ws.Cells.Item(350, 10).Value = "Mandarina"
$ws.Cells.Item(350, 11).Value = "Murcott"
$ws.Cells.Item(350, 12).Value = "Segunda"
$ws.Cells.Item(350, 13).Value = 50
$ws.Cells.Item(350, 14).Value = 7000
$ws.Cells.Item(350, 15).Value = 7000
$ws.Cells.Item(350, 16).Value = 7000
$ws.Cells.Item(350, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(350, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(350, 19).Value = 389
$ws.Cells.Item(350, 20).Value = 18
